# Append a new data row (row 19) to the extraction summary sheet, mirroring
# the existing rows' layout: A=library name (text), B-E=counts (numbers),
# F-I=percentages (stored as text, same as every other row), J=timestamp (text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "test"
$ws.Range("B19").Value = 1934
$ws.Range("C19").Value = 1661
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 1633

# Leading apostrophe forces these numeric-looking strings to be stored as
# text (matching the existing rows), then ClearFormats drops the resulting
# "number stored as text" quote-prefix style so no stray formatting is left
# behind on the new cells.
$ws.Range("F19").Value = "'36.96"
$ws.Range("G19").Value = "'31.75"
$ws.Range("H19").Value = "'0.08"
$ws.Range("I19").Value = "'31.21"
$ws.Range("F19:I19").ClearFormats()

$ws.Range("J19").Value = "2025-08-28 14:24:14"
